# [Feat] : CS_ATTACK, SC_POS_INTERPOLATION 프로토콜 수정
#
# 1) SC_POS_INTERPOLATION gains a new leading field: playerId (uint32, "플레이어 Id")
# 2) CS_ATTACK's old "bAttack" (bool) field becomes "hittedTargetId" (uint32, "피격 대상")
# 3) SC_ATTACK gains a new field right after playerId: hittedTargetId (uint32, "피격 대상")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) SC_POS_INTERPOLATION: insert playerId row right before the existing posX row (row 51) ---
$ws.Rows.Item(51).Insert()
$ws.Range("A51").Value = "SC_POS_INTERPOLATION"
$ws.Range("B51").Value = "playerId"
$ws.Range("C51").Value = "uint32"
$ws.Range("D51").Value = "플레이어 Id"

# --- 2) CS_ATTACK: the bAttack/bool row shifted down by the insert above, now at row 56 ---
$ws.Range("B56").Value = "hittedTargetId"
$ws.Range("C56").Value = "uint32"
$ws.Range("D56").Value = "피격 대상"

# --- 3) SC_ATTACK: insert hittedTargetId row right after playerId (old row 64, now row 65) ---
$ws.Rows.Item(65).Insert()
$ws.Range("A65").Value = "SC_ATTACK"
$ws.Range("B65").Value = "hittedTargetId"
$ws.Range("C65").Value = "uint32"
$ws.Range("D65").Value = "피격 대상"

# --- cosmetic: restore the view/selection reported in the saved file ---
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D77").Select()
